# Sync file from Google Drive
#
# This mirrors a re-sync of the "DDP_OUTPUT" bus-arrival snapshot workbook:
# the EstimatedTimeOfArrival timestamps (col F) and the derived
# MinutesToArrival counters (col O) were refreshed, and a handful of
# TypeOfBus (col L) / Load (col I) / Monitored (col J) values changed
# between the two Drive syncs, across all three NextBus* sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 (NextBus1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 45685.59071759259   # F2  EstimatedTimeOfArrival
$ws.Cells.Item(2, 15).Value = 2                  # O2  MinutesToArrival
$ws.Cells.Item(3, 15).Value = 2                  # O3  MinutesToArrival
$ws.Cells.Item(4, 6).Value = 45685.59165509259   # F4  EstimatedTimeOfArrival
$ws.Cells.Item(4, 15).Value = 3                  # O4  MinutesToArrival
$ws.Cells.Item(5, 6).Value = 45685.59210648148   # F5  EstimatedTimeOfArrival
$ws.Cells.Item(5, 12).Value = "BD"               # L5  TypeOfBus
$ws.Cells.Item(5, 15).Value = 4                  # O5  MinutesToArrival
$ws.Cells.Item(6, 6).Value = 45685.59135416667   # F6  EstimatedTimeOfArrival
$ws.Cells.Item(6, 15).Value = 3                  # O6  MinutesToArrival
$ws.Cells.Item(7, 6).Value = 45685.59237268518   # F7  EstimatedTimeOfArrival
$ws.Cells.Item(7, 15).Value = 4                  # O7  MinutesToArrival
$ws.Cells.Item(8, 6).Value = 45685.5909375       # F8  EstimatedTimeOfArrival
$ws.Cells.Item(8, 15).Value = 2                  # O8  MinutesToArrival
$ws.Cells.Item(9, 6).Value = 45685.59287037037   # F9  EstimatedTimeOfArrival
$ws.Cells.Item(9, 15).Value = 5                  # O9  MinutesToArrival
$ws.Cells.Item(10, 6).Value = 45685.59145833334  # F10 EstimatedTimeOfArrival
$ws.Cells.Item(11, 6).Value = 45685.59339120371  # F11 EstimatedTimeOfArrival
$ws.Cells.Item(11, 12).Value = "SD"              # L11 TypeOfBus
$ws.Cells.Item(11, 15).Value = 6                 # O11 MinutesToArrival
$ws.Cells.Item(12, 6).Value = 45685.59097222222  # F12 EstimatedTimeOfArrival
$ws.Cells.Item(12, 9).Value = "SEA"              # I12 Load
$ws.Cells.Item(12, 15).Value = 2                 # O12 MinutesToArrival
$ws.Cells.Item(13, 15).Value = 2                 # O13 MinutesToArrival
$ws.Cells.Item(14, 6).Value = 45685.59820601852  # F14 EstimatedTimeOfArrival
$ws.Cells.Item(14, 15).Value = 13                # O14 MinutesToArrival
$ws.Cells.Item(15, 6).Value = 45685.59876157407  # F15 EstimatedTimeOfArrival
$ws.Cells.Item(15, 15).Value = 13                # O15 MinutesToArrival

# --- Sheet 2 (NextBus2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 45685.60305555556   # F2  EstimatedTimeOfArrival
$ws.Cells.Item(2, 15).Value = 20                 # O2  MinutesToArrival
$ws.Cells.Item(3, 6).Value = 45685.59506944445   # F3  EstimatedTimeOfArrival
$ws.Cells.Item(3, 15).Value = 8                  # O3  MinutesToArrival
$ws.Cells.Item(4, 6).Value = 45685.59914351852   # F4  EstimatedTimeOfArrival
$ws.Cells.Item(5, 6).Value = 45685.59560185186   # F5  EstimatedTimeOfArrival
$ws.Cells.Item(5, 12).Value = "DD"               # L5  TypeOfBus
$ws.Cells.Item(5, 15).Value = 9                  # O5  MinutesToArrival
$ws.Cells.Item(6, 6).Value = 45685.59314814815   # F6  EstimatedTimeOfArrival
$ws.Cells.Item(6, 15).Value = 5                  # O6  MinutesToArrival
$ws.Cells.Item(7, 6).Value = 45685.60008101852   # F7  EstimatedTimeOfArrival
$ws.Cells.Item(7, 15).Value = 15                 # O7  MinutesToArrival
$ws.Cells.Item(8, 6).Value = 45685.59729166667   # F8  EstimatedTimeOfArrival
$ws.Cells.Item(8, 15).Value = 11                 # O8  MinutesToArrival
$ws.Cells.Item(9, 6).Value = 45685.59929398148   # F9  EstimatedTimeOfArrival
$ws.Cells.Item(9, 15).Value = 14                 # O9  MinutesToArrival
$ws.Cells.Item(10, 6).Value = 45685.59737268519  # F10 EstimatedTimeOfArrival
$ws.Cells.Item(10, 15).Value = 11                # O10 MinutesToArrival
$ws.Cells.Item(11, 6).Value = 45685.60056712963  # F11 EstimatedTimeOfArrival
$ws.Cells.Item(11, 12).Value = "DD"              # L11 TypeOfBus
$ws.Cells.Item(11, 15).Value = 16                # O11 MinutesToArrival
$ws.Cells.Item(12, 6).Value = 45685.60162037037  # F12 EstimatedTimeOfArrival
$ws.Cells.Item(13, 6).Value = 45685.60217592592  # F13 EstimatedTimeOfArrival
$ws.Cells.Item(13, 15).Value = 18                # O13 MinutesToArrival
$ws.Cells.Item(14, 6).Value = 45685.6047800926   # F14 EstimatedTimeOfArrival
$ws.Cells.Item(14, 15).Value = 22                # O14 MinutesToArrival
$ws.Cells.Item(15, 6).Value = 45685.60505787037  # F15 EstimatedTimeOfArrival
$ws.Cells.Item(15, 15).Value = 23                # O15 MinutesToArrival

# --- Sheet 3 (NextBus3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 45685.60755787037   # F2  EstimatedTimeOfArrival
$ws.Cells.Item(2, 15).Value = 26                 # O2  MinutesToArrival
$ws.Cells.Item(3, 6).Value = 45685.59855324074   # F3  EstimatedTimeOfArrival
$ws.Cells.Item(3, 15).Value = 13                 # O3  MinutesToArrival
$ws.Cells.Item(4, 6).Value = 45685.60621527778   # F4  EstimatedTimeOfArrival
$ws.Cells.Item(4, 15).Value = 24                 # O4  MinutesToArrival
$ws.Cells.Item(5, 6).Value = 45685.60087962963   # F5  EstimatedTimeOfArrival
$ws.Cells.Item(5, 12).Value = "SD"               # L5  TypeOfBus
$ws.Cells.Item(5, 15).Value = 17                 # O5  MinutesToArrival
$ws.Cells.Item(6, 6).Value = 45685.60180555555   # F6  EstimatedTimeOfArrival
$ws.Cells.Item(6, 15).Value = 18                 # O6  MinutesToArrival
$ws.Cells.Item(7, 15).Value = 24                 # O7  MinutesToArrival
$ws.Cells.Item(8, 6).Value = 45685.60627314815   # F8  EstimatedTimeOfArrival
$ws.Cells.Item(8, 15).Value = 24                 # O8  MinutesToArrival
$ws.Cells.Item(9, 6).Value = 45685.60402777778   # F9  EstimatedTimeOfArrival
$ws.Cells.Item(9, 15).Value = 21                 # O9  MinutesToArrival
$ws.Cells.Item(10, 6).Value = 45685.60359953704  # F10 EstimatedTimeOfArrival
$ws.Cells.Item(10, 15).Value = 20                # O10 MinutesToArrival
$ws.Cells.Item(11, 6).Value = 45685.60436342593  # F11 EstimatedTimeOfArrival
$ws.Cells.Item(11, 12).Value = "SD"              # L11 TypeOfBus
$ws.Cells.Item(11, 15).Value = 22                # O11 MinutesToArrival
$ws.Cells.Item(12, 15).Value = 30                # O12 MinutesToArrival
$ws.Cells.Item(13, 6).Value = 45685.61067129629  # F13 EstimatedTimeOfArrival
$ws.Cells.Item(13, 15).Value = 31                # O13 MinutesToArrival
$ws.Cells.Item(14, 6).Value = 45685.6133912037   # F14 EstimatedTimeOfArrival
$ws.Cells.Item(14, 15).Value = 35                # O14 MinutesToArrival
$ws.Cells.Item(15, 6).Value = 45685.61273148148  # F15 EstimatedTimeOfArrival
$ws.Cells.Item(15, 10).Value = 0                 # J15 Monitored
$ws.Cells.Item(15, 15).Value = 34                # O15 MinutesToArrival
